# Add kernel-SVR parameters (columns K:M) to the parameter table on Sheet1
# and move the current selection, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1) entries for the added parameters.
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# New values row (row 2) for the added parameters.
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Best-effort: reflect the author's window-size/position bookkeeping.
# (Cosmetic workbook-view metadata; harmless if the host doesn't persist it.)
$win = $excel.ActiveWindow
$win.Left = -98
$win.Top = -98
$win.Width = 23236
$win.Height = 13875

# Move the active selection to K7, as recorded in the saved file.
$ws.Range("K7").Select()
